$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G10").Value = "a"
$ws.Range("G11").Value = "b"
$ws.Range("G12").Value = "c"
$ws.Range("G13").Value = "d"
$ws.Range("G14").Value = "e"
$ws.Range("G15").Value = "f"
